$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.757.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.293.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '96.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.52%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.92'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.637.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.850'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.290.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.707.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.69%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +11.90%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.63%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +8.15%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.97'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.88'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.45%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.43'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.60%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0352'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.105'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.237'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.37%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.33%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.41%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Aave'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +13.52%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.433'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.82%  '
$ws.Range('E51').Style = 'Normal'
